$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change A3 value from "Adam" to "Madam"
$ws.Range("A3").Value = "Madam"

# Update the active selection to A3
$ws.Range("A3").Select()
